$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Stores")

# Unmerge the A1:F2 cell range
$ws.Range("A1:F2").UnMerge()

# Remove the center alignment from A1 (and the other previously-merged cells)
$ws.Range("A1:F2").HorizontalAlignment = 1

# Add new data row: "Number of items" = 9
$ws.Range("B4").Value = "Number of items"
$ws.Range("C4").Value = 9

# Update selected cell to D6
$ws.Range("D6").Select()

$wb.Save()
